# Daily attendance processing - 2025-11-01 10:51:32
# Reorders the "Recorded By" (column G) values for the "backup@backdoor.com"
# and "dnasr281@gmail.com" accounts so the "System"/"system" token(s) come
# first in the comma-separated list, with the email identifier moved to the
# end (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# The "admin@admin.com" rows are intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Exact before -> after mapping, matching only the affected "Recorded By"
# values (leaving every other value, e.g. the admin@admin.com rows, alone).
$map = @{
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "backup@backdoor.com, System, system" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ([string]::IsNullOrEmpty($value)) { continue }

    if ($map.ContainsKey($value)) {
        $cell.Value = $map[$value]
    }
}
